$wb = $excel.ActiveWorkbook

# Helper: set a cell to a literal TEXT value even if it looks numeric,
# without leaving a visible style on the cell.
function Set-TextCell($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# ---- Add CaseDetailStat sheet (copy of StatOutput) ----
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "CaseDetailStat"
$newSheet.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

$ws6 = $wb.Worksheets.Item("CaseDetailStat")

$ws6.Range("A1").Value = 'number_of_files'
$ws6.Range("B1").Value = 'number_of_sample'
$ws6.Range("C1").Value = 'number_of_cases'
$ws6.Range("D1").Value = 'number_of_study'

Set-TextCell $ws6.Range("A2") '331'
Set-TextCell $ws6.Range("B2") '136'
Set-TextCell $ws6.Range("C2") '60'
Set-TextCell $ws6.Range("D2") '1'

# ---- Add CaseDetailStat_Message sheet ----
$newSheet2 = $wb.Worksheets.Add()
$newSheet2.Name = "CaseDetailStat_Message"
$newSheet2.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

$ws7 = $wb.Worksheets.Item("CaseDetailStat_Message")

$ws7.Range("A1").Value = 'Neo4j_URL:'
$ws7.Range("A2").Value = 'bolt://ncias-q2251-c.nci.nih.gov:7687'
$ws7.Range("A3").Value = 'User_name:'
$ws7.Range("A4").Value = 'neo4j'
$ws7.Range("A5").Value = 'PWD:'
$ws7.Range("A6").Value = 'icdcDBneo4j0'
$ws7.Range("A7").Value = 'Cypher:'
$ws7.Range("A8").Value = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE s.clinical_study_designation IN [''NCATS-COP01''] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(s.clinical_study_designation,'''') AS `Study Code` , coalesce(s.clinical_study_type,'''') AS  `Study Type`, coalesce(demo.breed,'''') AS Breed , coalesce(diag.disease_term,'''') AS Diagnosis , coalesce(diag.stage_of_disease,'''') AS `Stage of Disease` ,  coalesce(demo.patient_age_at_enrollment,'''') AS Age , coalesce(demo.sex,'''') AS Sex , coalesce(demo.neutered_indicator,'''') AS  `Neutered Status`'
$ws7.Range("A9").Value = 'Output:'
$ws7.Range("A10").Value = 'C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC02_Canine_Filter_Study-NCATS_Neo4jData.xlsx'

$ws7.Range("A11").Value = 'Neo4j_URL:'
$ws7.Range("A12").Value = 'bolt://ncias-q2251-c.nci.nih.gov:7687'
$ws7.Range("A13").Value = 'User_name:'
$ws7.Range("A14").Value = 'neo4j'
$ws7.Range("A15").Value = 'PWD:'
$ws7.Range("A16").Value = 'icdcDBneo4j0'
$ws7.Range("A17").Value = 'Cypher:'
$ws7.Range("A18").Value = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE s.clinical_study_designation IN [''NCATS-COP01'']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study'
$ws7.Range("A19").Value = 'Output:'
$ws7.Range("A20").Value = 'C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC02_Canine_Filter_Study-NCATS_Neo4jData.xlsx'

$ws7.Range("A21").Value = 'Neo4j_URL:'
$ws7.Range("A22").Value = 'bolt://ncias-q2251-c.nci.nih.gov:7687'
$ws7.Range("A23").Value = 'User_name:'
$ws7.Range("A24").Value = 'neo4j'
$ws7.Range("A25").Value = 'PWD:'
$ws7.Range("A26").Value = 'icdcDBneo4j0'
$ws7.Range("A27").Value = 'Cypher:'
$ws7.Range("A28").Value = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE s.clinical_study_designation IN [''NCATS-COP01'']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study'
$ws7.Range("A29").Value = 'Output:'
$ws7.Range("A30").Value = 'C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC02_Canine_Filter_Study-NCATS_Neo4jData.xlsx'

